$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Determined which features to implement: assign "Annabelle" to
# "Create Tutorials" and "Add invisible watermark in metadata"
$ws.Range("C9").Value = "Annabelle"
$ws.Range("C10").Value = "Annabelle"

$ws.Range("C9").Select()
